# Update "想去人数" (column F) counts across the four worksheets.
# Mapping of Sheet name -> { Row -> NewValue } built from the source diff.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        5  = 2310
        7  = 8181
        13 = 4518
        15 = 798
        20 = 6533
        24 = 4403
        25 = 321
        27 = 2059
        28 = 1198
        31 = 74
        33 = 47
        34 = 92
        35 = 337
        38 = 153
        41 = 1233
        43 = 73
        48 = 35
        49 = 25
    }
    "演出" = @{
        6  = 698
        7  = 405
        10 = 221
        17 = 111
        26 = 179
        30 = 13
    }
    "本地生活" = @{
        4  = 469
        8  = 3165
        9  = 1088
        11 = 1529
        12 = 1869
        13 = 357
        14 = 224
    }
    "全部类型" = @{
        2  = 469
        7  = 3165
        8  = 2310
        9  = 1088
        12 = 1529
        14 = 698
        16 = 1869
        17 = 4518
        18 = 357
        19 = 405
        21 = 798
        26 = 6533
        28 = 224
        30 = 4403
        31 = 321
        32 = 2059
        33 = 1198
        36 = 74
        37 = 47
        39 = 92
        40 = 337
        42 = 153
        44 = 1233
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
